$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header row copy (typos -> proper header names), add two new headers ---
$ws.Range("B1").Value = "Employee ID"
$ws.Range("C1").Value = "Designation"
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "Department"
$ws.Range("F1").Value = "Manager Name"
$ws.Range("G1").Value = "Manager Email"
$ws.Range("H1").Value = "Joining Date"
$ws.Range("I1").Value = ""

# Match the formatting (style index) already used across row 1 for the newly
# added H1/I1 header cells.
$ws.Range("G1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)

# --- Add a sample data row ---
$ws.Range("A2").Value = "ABCD"
$ws.Range("B2").Value = 12322
$ws.Range("C2").Value = "Software Engineer"
$ws.Range("D2").Value = "abc@gmail.com"
$ws.Range("E2").Value = "Development"
$ws.Range("F2").Value = "XYZ"
$ws.Range("G2").Value = "xyz@gmail.com"

# Give the new data row the same formatting as the header row cells.
$ws.Range("G1").Copy()
$ws.Range("A2:G2").PasteSpecial(-4122)

# Joining date, formatted as mm/dd/yyyy (set the number format first so the
# new style picks up the right font / doesn't get reset when assigning the
# date value).
$ws.Range("H2").NumberFormat = "mm/dd/yyyy"
$ws.Range("H2").Value = (Get-Date -Year 2024 -Month 1 -Day 1).Date
